$wb = $excel.ActiveWorkbook

# --- Add the new "plt_ref" worksheet at the end of the workbook ---
$newWs = $wb.Worksheets.Add()
$newWs.Name = "plt_ref"
$newWs.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

$ws4 = $wb.Worksheets.Item("plt_ref")

$ws4.Range("A1").Value = "# other useful name referencing based on project specific naming"
$ws4.Range("A2").Value = "flux"
$ws4.Range("A3").Value = "Qbase"
$ws4.Range("A4").Value = "Qrech"
$ws4.Range("A5").Value = "flowing"
$ws4.Range("A6").Value = "Qout"
$ws4.Range("B2").Value = "color"
$ws4.Range("C2").Value = "name"
$ws4.Range("C3").Value = "Baseflow"
$ws4.Range("C4").Value = "Stream Losses"
$ws4.Range("C5").Value = "No. Days with Flow"
$ws4.Range("C6").Value = "Streamflow"
$ws4.Range("A7").Value = "num_sfr_coarse"
$ws4.Range("A8").Value = "num_lak_coarse"
$ws4.Range("C7").Value = "No. of Coarse Reaches"
$ws4.Range("C8").Value = "No. of Coarse Lake Cells"

# --- Update selection on owhm_wb_dict ---
$ws2 = $wb.Worksheets.Item("owhm_wb_dict")
$ws2.Activate()
$ws2.Range("A2:C2").Select()

# --- Make plt_ref the active/selected sheet with C9 selected ---
$ws4.Activate()
$ws4.Range("C9").Select()

Write-Host "done"
